$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row above the current row 133, shifting the
# existing rows 133-208 down to 134-209 (dimension grows to A1:R209).
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new record.
$ws.Range("A133").Value = 9
$ws.Range("B133").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C133").Value = "Metropolitana"
$ws.Range("D133").Value = 44455
$ws.Range("E133").Value = 13
$ws.Range("F133").Value = 100112044
$ws.Range("G133").Value = "Perejil"
$ws.Range("H133").Value = "Sin especificar"
$ws.Range("I133").Value = "Primera"
$ws.Range("J133").Value = 106
$ws.Range("K133").Value = 8000
$ws.Range("L133").Value = 10000
$ws.Range("M133").Value = 9000
$ws.Range("N133").Value = "`$/docena de atados"
$ws.Range("O133").Value = "Región Metropolitana"
$ws.Range("P133").Value = 3000
$ws.Range("Q133").Value = 3
$ws.Range("R133").Value = "Hortaliza"
